# GPLIM-2588 Fix spreadsheet headers.
# A1 was "Sample ID" -> becomes "Specimen_Number"
# F1 was "T/N" -> becomes "SAMPLE_TYPE"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe preserves the existing cell style's quote-prefix flag
# (style used by A1 already has quotePrefix set) instead of Excel cloning a
# brand new style without it.
$ws.Range("A1").Value = "'Specimen_Number"
$ws.Range("F1").Value = "SAMPLE_TYPE"

$ws.Range("A2").Select()
